$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.408.75"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.848.40"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.36"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6261"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07497"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2903"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.43"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.848.54"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.003"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6808"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001045"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.26"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "2.103.79"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.172"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "29.438.28"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.85"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.476"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9991"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.10"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1374"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.413"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.55"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06401"
$ws.Range("E29").Value = "  +14.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.477"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.094"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.094"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6986"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.577"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "1.265.14"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.825"
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01830"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.629"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9101"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "2.009.14"
$ws.Range("E44").Value = "  -18.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.48"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.30"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.745"
$ws.Range("E47").Value = "  +4.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.076"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1173"
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000117"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.043"
$ws.Range("E51").Value = "  +0.58%  "
